# Handles float input without breaking stuff
#
# The quiz marksheet previously had three "Student Ans / Correct Ans" blocks
# (A:B, D:E, G:H). The grading for this student only produced two real
# blocks of results, so:
#   - the third block (columns G:H) is removed entirely,
#   - the second block (D:E) is trimmed down to only the rows that still
#     have data (rows 16-18),
#   - the first block (A:B) gets the student's answers filled into column A,
#     using the "correct" (green) style for matches and leaving a blank
#     "not attempted" (black) style otherwise,
#   - the summary rows (10-12) are updated with the recomputed counts/score.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Drop the unused third answer block (columns G:H) completely ---
$ws.Range("G1:H100").Clear() | Out-Null

# --- Drop the now-empty tail rows of the second answer block (D19:E40) ---
$ws.Range("D19:E40").Clear() | Out-Null

# --- Re-style the "No./Marking/Total" row labels like the header row (A9) ---
foreach ($cellRef in @("A10", "A11", "A12")) {
    $ws.Range("A9").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats) | Out-Null
}
$ws.Range("A10").Value = "No."
$ws.Range("A11").Value = "Marking"
$ws.Range("A12").Value = "Total"

# --- Fill in student answers that match the correct answer (green "correct" style, copied from B10) ---
$correctCells = @{
    "A16" = "Option A"; "D16" = "Option A"
    "A17" = "Option D"
    "A18" = "Option B"; "D18" = "Option D"
    "A19" = "Option C"
    "A20" = "Option B"
    "A21" = "Option C"
    "A22" = "Option D"
    "A24" = "Option A"
    "A25" = "Option A"
    "A28" = "Option D"
    "A31" = "Option D"
    "A32" = "Option C"
    "A33" = "Option D"
    "A34" = "Option B"
    "A35" = "Option D"
    "A36" = "Option A"
    "A37" = "Option A"
    "A38" = "Option A"
    "A39" = "Option D"
    "A40" = "Option D"
}
foreach ($cellRef in $correctCells.Keys) {
    $ws.Range("B10").Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial($xlPasteFormats) | Out-Null
}
foreach ($cellRef in $correctCells.Keys) {
    $ws.Range($cellRef).Value = $correctCells[$cellRef]
}

# --- D17 is a student answer that is wrong (red "incorrect" style, copied from C10) ---
$ws.Range("C10").Copy() | Out-Null
$ws.Range("D17").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D17").Value = "Option D"

# --- Recomputed summary numbers (row 10: No., row 11: Marking, row 12: Total) ---
$ws.Range("B10").Value = 22
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 88
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "87/112"
